$d = $word.ActiveDocument

# Apply edits from the end of the document towards the start so that
# earlier (lower) character offsets remain valid while later edits
# change the document length.

# Hunk 6: <meta name="url" ...> - point at the 2025 article instead of the 2024 one.
$rng = $d.Range(625, 708)
$rng.Text = "Enlightenment/Articles/2025/1-Blender-Continued/7-Loop-Tools/8-Space/8-Space.html"

# Hunk 5: <meta name="revised" ...> - bump the revision date.
$rng = $d.Range(534, 562)
$rng.Text = "Monday, February 24, 2025"

# Hunk 4: <meta name="category" ...> - insert "The Space Loop Tool, " before the
# existing "Blender, 3D Modeling, Animation, Graphic Art" keyword run, leaving the
# surrounding runs untouched.
$rng = $d.Range(368, 368)
$rng.InsertBefore("The Space Loop Tool, ")

# Hunk 3: <meta name="description" ...> - rewrite the explanatory sentence.
$rng = $d.Range(280, 335)
$rng.Text = "how to use the Space Tool inside of the Loop Tool" + [char]8217 + "s collection/>"

# Hunk 2: <meta name="keywords" ...> - replace the leading blank placeholder with real text.
$rng = $d.Range(172, 173)
$rng.Text = "The Space Loop Tool, "

# Hunk 1: <meta name="title" ...> - new article title.
$rng = $d.Range(28, 45)
$rng.Text = "8 Space"

Write-Host "edits applied"
